# Added medium-complexity steel model
#
# 1. Fills in "air" as the Carbon Fate (col C) for the EUROFER rows that
#    previously had no value there.
# 2. Appends 16 new base-case scenario rows (global/EU/USA/China/India/
#    Japan/Russia/SouthKorea, each split into a BF+BOF and an EAF route)
#    describing the new medium-complexity steel model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "air" for Carbon Fate (col C) on the EUROFER rows that were
# --- missing it (rows 11 and 15 already had a different Carbon Fate value).
$airRows = @(7, 8, 9, 10, 12, 13, 14, 16)
foreach ($r in $airRows) {
    $ws.Cells.Item($r, 3).Value = "air"
}

# Row -> Scenario name (col A), in the order the names were typed in
# (grouped by region, BF row immediately followed by its EAF row) so the
# shared-string table comes out in the same sequence as the source edit.
$aOrder = @(
    @{ Row = 18; Name = "global-BF-base" },
    @{ Row = 26; Name = "global-EAF-base" },
    @{ Row = 20; Name = "EU-BF-base" },
    @{ Row = 28; Name = "EU-EAF-base" },
    @{ Row = 25; Name = "USA-BF-base" },
    @{ Row = 33; Name = "USA-EAF-base" },
    @{ Row = 19; Name = "China-BF-base" },
    @{ Row = 27; Name = "China-EAF-base" },
    @{ Row = 21; Name = "India-BF-base" },
    @{ Row = 22; Name = "Japan-BF-base" },
    @{ Row = 23; Name = "Russia-BF-base" },
    @{ Row = 24; Name = "SouthKorea-BF-base" },
    @{ Row = 29; Name = "India-EAF-base" },
    @{ Row = 30; Name = "Japan-EAF-base" },
    @{ Row = 31; Name = "Russia-EAF-base" },
    @{ Row = 32; Name = "SouthKorea-EAF-base" }
)
foreach ($row in $aOrder) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Name
}

# Row -> Carbon Sources (col B), filled row by row.
$bValues = @{
    18 = "coal"; 19 = "coal"; 20 = "coal"; 21 = "coal";
    22 = "coal, waste plastics"; 23 = "coal, natural gas";
    24 = "coal"; 25 = "coal";
    26 = "electricity mix"; 27 = "electricity mix"; 28 = "electricity mix";
    29 = "electricity mix"; 30 = "electricity mix"; 31 = "electricity mix";
    32 = "electricity mix"; 33 = "electricity mix"
}
foreach ($r in 18..33) {
    $ws.Cells.Item($r, 2).Value = $bValues[$r]
}

# Row -> Carbon Fate (col C) - "air" for every new row.
foreach ($r in 18..33) {
    $ws.Cells.Item($r, 3).Value = "air"
}

# Row -> Technology Notes (col D) - BF+BOF for the coal-based rows, EAF for
# the electricity-based rows.
$dValues = @{
    18 = "BF+BOF"; 19 = "BF+BOF"; 20 = "BF+BOF"; 21 = "BF+BOF";
    22 = "BF+BOF"; 23 = "BF+BOF"; 24 = "BF+BOF"; 25 = "BF+BOF";
    26 = "EAF"; 27 = "EAF"; 28 = "EAF"; 29 = "EAF";
    30 = "EAF"; 31 = "EAF"; 32 = "EAF"; 33 = "EAF"
}
foreach ($r in 18..33) {
    $ws.Cells.Item($r, 4).Value = $dValues[$r]
    $ws.Cells.Item($r, 4).WrapText = $true
}

# --- Re-freeze the header row / first column, matching the refreshed view.
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D36").Select()
